# Insert a new row of data at row 25 (pushing existing rows 25-59 down to 26-60)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 25; this shifts rows 25..59 -> 26..60
# and carries forward the existing date-format styling on column D, matching the
# workbook's observed "insert row" behavior.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record's values.
$ws.Cells.Item(25, 1).Value = 11
$ws.Cells.Item(25, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(25, 3).Value = "Bíobío"
$ws.Cells.Item(25, 4).Value = 44778
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(25, 6).Value = 100112013
$ws.Cells.Item(25, 7).Value = "Alcachofa"
$ws.Cells.Item(25, 8).Value = "Argentina(o)"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 150
$ws.Cells.Item(25, 11).Value = 14000
$ws.Cells.Item(25, 12).Value = 15000
$ws.Cells.Item(25, 13).Value = 14533
$ws.Cells.Item(25, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(25, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 16).Value = 291
$ws.Cells.Item(25, 17).Value = 50
$ws.Cells.Item(25, 18).Value = "Hortaliza"
